$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.805.53"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "3.341.68"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.30"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "655.86"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.41"
$ws.Range("E7").Value = "  -4.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.423"
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -5.49%  "
$ws.Range("D11").Value = "3.339.96"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.207"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.71"
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("D14").Value = "96.556.03"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.10"
$ws.Range("E15").Value = "  -3.35%  "
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "3.967.82"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.72"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").Value = "3.359.80"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.574"
$ws.Range("E20").Value = "  +14.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.69"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "508.69"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.34"
$ws.Range("E24").Value = "  -2.92%  "
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.62"
$ws.Range("E26").Value = "  +6.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.56"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("E28").Value = "  -4.45%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.517.87"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.146"
$ws.Range("E30").Value = "  -4.18%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.38"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.189"
$ws.Range("E33").Value = "  -6.51%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.53"
$ws.Range("E34").Value = "  +11.78%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.554"
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "28.37"
$ws.Range("E37").Value = "  -4.86%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.50"
$ws.Range("E38").Value = "  +4.98%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.84"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.151"
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "507.28"
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.35"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0433"
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.837"
$ws.Range("E45").Value = "  -4.54%  "
$ws.Range("B46").Value = "MantraDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.70"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.68"
$ws.Range("E47").Value = "  +5.97%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.58"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.51"
$ws.Range("E49").Value = "  +2.93%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.67"
$ws.Range("E50").Value = "  +6.67%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.11"
$ws.Range("E51").Value = "  -5.65%  "
